# Apply the edits described by the commit "adding formulas to Excel output"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet/tab: "Sensed Comp3_LR History" -> "Sensed Comp3_LR"
$ws.Name = "Sensed Comp3_LR"

# 2. Give columns A:D an explicit custom width (~12.29 characters in the
#    underlying OOXML "width" units). The COM ColumnWidth property is offset
#    by ~0.8333 characters from the raw OOXML width and is itself quantized
#    to the nearest 1/6 of a character by this engine, so we back-solve for
#    the COM value that lands closest to the target raw width of 12.28515625.
$ws.Range("A1:D1").EntireColumn.ColumnWidth = 11.451822916666666

# 3. Turn A1 from the shared string "Time Step" into a formula cell
$ws.Range("A1").Formula = "=MODE(C[100, 0, 0, 3] : C[100, 0, 0, 3])"

# 4. Update the sensor-reading cells that changed
$ws.Range("C10:C13").Value = 0
$ws.Range("B20:B21").Value = 2
$ws.Range("B22:B41").Value = 0
$ws.Range("C60:C63").Value = 0
$ws.Range("B70:B71").Value = 2
$ws.Range("B72:B91").Value = 0
